$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the new CTS cancellation record.
$ws.Range("C2").Value = "'1001615285"
$ws.Range("F2").ClearContents()
$ws.Range("G2").Value = "AASIM23200V88KK9Q"
$ws.Range("H2").ClearContents()

# Update the saved selection on the sheet.
$ws.Range("E6").Select()
